$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text-formatted so numeric-looking
# strings like "302.17" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.057.53"
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").Value = "1.598.68"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  -0.34%  "

$ws.Range("E5").Value = "  -0.25%  "

$ws.Range("D6").Value = "302.17"
$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("E7").Value = "  +0.37%  "

$ws.Range("D8").Value = "0.3624"
$ws.Range("E8").Value = "  -1.13%  "

$ws.Range("D9").Value = "51.13"
$ws.Range("E9").Value = "  +6.56%  "

$ws.Range("D10").Value = "1.241"
$ws.Range("E10").Value = "  -2.77%  "

$ws.Range("E11").Value = "  -0.35%  "

$ws.Range("D12").Value = "0.08119"
$ws.Range("E12").Value = "  +0.57%  "

$ws.Range("D13").Value = "22.19"
$ws.Range("E13").Value = "  -3.03%  "

$ws.Range("D14").Value = "6.549"
$ws.Range("E14").Value = "  -1.41%  "

$ws.Range("D15").Value = "7.304"
$ws.Range("E15").Value = "  -3.47%  "

$ws.Range("D16").Value = "0.00001236"
$ws.Range("E16").Value = "  -2.07%  "

$ws.Range("D17").Value = "1.595.59"
$ws.Range("E17").Value = "  -0.06%  "

$ws.Range("D18").Value = "92.32"
$ws.Range("E18").Value = "  +0.37%  "

$ws.Range("D19").Value = "0.06847"
$ws.Range("E19").Value = "  +0.46%  "

$ws.Range("D20").Value = "18.11"
$ws.Range("E20").Value = "  -1.90%  "

$ws.Range("D21").Value = "6.492"
$ws.Range("E21").Value = "  -1.73%  "

$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").Value = "12.95"
$ws.Range("E23").Value = "  -1.13%  "

$ws.Range("D24").Value = "23.063.16"
$ws.Range("E24").Value = "  -0.34%  "

$ws.Range("D25").Value = "2.376"
$ws.Range("E25").Value = "  +0.63%  "

$ws.Range("D26").Value = "2.797"
$ws.Range("E26").Value = "  -4.99%  "

$ws.Range("D27").Value = "21.04"
$ws.Range("E27").Value = "  -0.60%  "

$ws.Range("D28").Value = "149.04"
$ws.Range("E28").Value = "  -1.41%  "

$ws.Range("D29").Value = "5.250"
$ws.Range("E29").Value = "  +0.25%  "

$ws.Range("D30").Value = "133.90"
$ws.Range("E30").Value = "  +0.99%  "

$ws.Range("D31").Value = "2.380"
$ws.Range("E31").Value = "  -3.06%  "

$ws.Range("D32").Value = "6.720"
$ws.Range("E32").Value = "  -6.44%  "

$ws.Range("D33").Value = "1.772.96"
$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("D34").Value = "0.9574"
$ws.Range("E34").Value = "  -1.64%  "

$ws.Range("D35").Value = "0.07486"
$ws.Range("E35").Value = "  -3.37%  "

$ws.Range("D36").Value = "0.02701"
$ws.Range("E36").Value = "  -2.74%  "

$ws.Range("D37").Value = "10.14"
$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("D38").Value = "6.163"
$ws.Range("E38").Value = "  -2.01%  "

$ws.Range("D39").Value = "0.2506"
$ws.Range("E39").Value = "  -1.46%  "

$ws.Range("D40").Value = "0.08809"
$ws.Range("E40").Value = "  -0.53%  "

$ws.Range("D41").Value = "1.360"
$ws.Range("E41").Value = "  -2.12%  "

$ws.Range("D42").Value = "0.7007"
$ws.Range("E42").Value = "  -2.15%  "

$ws.Range("D43").Value = "12.39"
$ws.Range("E43").Value = "  -2.75%  "

$ws.Range("D44").Value = "15.09"
$ws.Range("E44").Value = "  -6.38%  "

$ws.Range("D45").Value = "0.6532"
$ws.Range("E45").Value = "  -1.36%  "

$ws.Range("D46").Value = "4.010"
$ws.Range("E46").Value = "  +1.11%  "

$ws.Range("E47").Value = "  -1.76%  "

$ws.Range("D48").Value = "132.13"
$ws.Range("E48").Value = "  +1.01%  "

$ws.Range("D49").Value = "0.07919"
$ws.Range("E49").Value = "  -0.86%  "

$ws.Range("D50").Value = "1.214"
$ws.Range("E50").Value = "  +3.86%  "

$ws.Range("D51").Value = "1.225"
$ws.Range("E51").Value = "  +4.07%  "
